$d = $word.ActiveDocument
$d.Content.Find.Execute("[.]Net 4.5.1", $false, $false, $true, $false, $false, `
                         $true, 1, $false, ".Net 4.5.2", 2)
